# Appends ten new paragraphs at the very end of the document body (right
# after the last, already-empty, paragraph and before the sectPr):
#   - 4 plain empty paragraphs
#   - 5 bold paragraphs of notes text (one of them with an internal line
#     break) and 1 empty-but-bold paragraph in between them, matching the
#     target OOXML exactly, including <w:b/><w:bCs/> run/paragraph-mark
#     formatting.
#
# Word's Range.InsertXML only *replaces* the content of the range it is
# called on, so we first mint a fresh empty paragraph at the tail of the
# story (InsertParagraphAfter) and target that whole paragraph (its mark
# included) with one InsertXML call carrying every new paragraph as a
# WordprocessingML package fragment. That reliably reproduces <w:b/> and
# <w:bCs/> (which the plain Range.Bold / Font.Bold COM properties do not
# expose). InsertXML'ing a paragraph-ended fragment always leaves one
# extra trailing empty paragraph behind (because the final paragraph mark
# of the fragment pushes the story's own end-of-range mark into a new
# paragraph) -- that spare paragraph is removed at the end.

$d = $word.ActiveDocument

$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter() | Out-Null
$target = $d.Paragraphs.Last

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> ADD D0 RESP VS NON-RESP</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">RUN AGAIN WITH DATA CLEANED WITHOUT NORMALISATION </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&#8211;</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Marie sent</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Run with SPE and Insoluble data (combine into single dataset for logistic?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Further comparisons:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>D0 vs other timepoints</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml) | Out-Null

$secondLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$trailing = $d.Paragraphs.Last
$spare = $d.Range($secondLast.Range.End - 1, $trailing.Range.End)
$spare.Delete() | Out-Null

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
